$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 4432.346
$ws.Range("J53").Value = 9307.75
$ws.Range("L53").Value = 9307.75
$ws.Range("N53").Value = -10581.75
$ws.Range("H76").Value = 6950282
$ws.Range("J76").Value = 6974.75
$ws.Range("L76").Value = 6974.75
$ws.Range("N76").Value = -7604.75
$ws.Range("H79").Value = 6950282
$ws.Range("J79").Value = 6974.75
$ws.Range("L79").Value = 6974.75
$ws.Range("N79").Value = -9158.75
$ws.Range("H80").Value = 645.4
$ws.Range("I80").Value = 620.4286
$ws.Range("K80").Value = 1861.2858
$ws.Range("M80").Value = -863.2857999999999
$ws.Range("H83").Value = 645.4
$ws.Range("I83").Value = 620.4286
$ws.Range("K83").Value = 5583.8574
$ws.Range("M83").Value = -591.8573999999999
$ws.Range("H86").Value = 3049.8333
$ws.Range("I86").Value = 2850
$ws.Range("J86").Value = 3449.5
$ws.Range("K86").Value = 2850
$ws.Range("L86").Value = 3449.5
$ws.Range("M86").Value = -1727
$ws.Range("N86").Value = -5695.5
$ws.Range("H89").Value = 3049.8333
$ws.Range("I89").Value = 2850
$ws.Range("J89").Value = 3449.5
$ws.Range("K89").Value = 14250
$ws.Range("L89").Value = 17247.5
$ws.Range("M89").Value = -8634
$ws.Range("N89").Value = -28479.5
$ws.Range("H115").Value = 313.5
$ws.Range("I115").Value = 313.5
$ws.Range("K115").Value = 940.5
$ws.Range("M115").Value = 626.5
$ws.Range("H121").Value = 2218
$ws.Range("J121").Value = 2218
$ws.Range("L121").Value = 6654
$ws.Range("N121").Value = -10148
$ws.Range("H138").Value = 2665.4905
$ws.Range("I138").Value = 1182.7097
$ws.Range("J138").Value = 4754.864
$ws.Range("K138").Value = 3548.1291
$ws.Range("L138").Value = 14264.592
$ws.Range("M138").Value = 1591.8709
$ws.Range("N138").Value = -24544.592

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1199.1875
$ws.Range("I2").Value = 830.375
$ws.Range("J2").Value = 1568
$ws.Range("K2").Value = 830.375
$ws.Range("L2").Value = 1568
$ws.Range("M2").Value = -717.375
$ws.Range("N2").Value = -1794
$ws.Range("H74").Value = 71654
$ws.Range("I74").Value = 7541.52
$ws.Range("K74").Value = 7541.52
$ws.Range("M74").Value = -6667.52
$ws.Range("H77").Value = 71654
$ws.Range("I77").Value = 7541.52
$ws.Range("K77").Value = 37707.60000000001
$ws.Range("M77").Value = -33339.60000000001
$ws.Range("H92").Value = 79999.5
$ws.Range("J92").Value = 79999.5
$ws.Range("L92").Value = 79999.5
$ws.Range("N92").Value = -84991.5
$ws.Range("H116").Value = 1199.1875
$ws.Range("I116").Value = 830.375
$ws.Range("J116").Value = 1568
$ws.Range("K116").Value = 830.375
$ws.Range("L116").Value = 1568
$ws.Range("M116").Value = 1463.625
$ws.Range("N116").Value = -6156
$ws.Range("H132").Value = 2774.7
$ws.Range("I132").Value = 2550.6667
$ws.Range("J132").Value = 3110.75
$ws.Range("K132").Value = 7652.000100000001
$ws.Range("L132").Value = 9332.25
$ws.Range("M132").Value = -5122.000100000001
$ws.Range("N132").Value = -14392.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1199.1875
$ws.Range("I3").Value = 830.375
$ws.Range("J3").Value = 1568
$ws.Range("K3").Value = 830.375
$ws.Range("L3").Value = 1568
$ws.Range("M3").Value = -716.375
$ws.Range("N3").Value = -1796
$ws.Range("H82").Value = 27328.715
$ws.Range("I82").Value = 10260.4
$ws.Range("K82").Value = 10260.4
$ws.Range("M82").Value = -9877.4
$ws.Range("H85").Value = 27328.715
$ws.Range("I85").Value = 10260.4
$ws.Range("K85").Value = 10260.4
$ws.Range("M85").Value = -8934.4
$ws.Range("H96").Value = 17899
$ws.Range("I96").Value = 17899
$ws.Range("K96").Value = 17899
$ws.Range("M96").Value = -15153
$ws.Range("H134").Value = 3273.7646
$ws.Range("I134").Value = 1655.258
$ws.Range("K134").Value = 4965.774
$ws.Range("M134").Value = -2430.774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4269.53
$ws.Range("I31").Value = 1374.5625
$ws.Range("J31").Value = 4820.952
$ws.Range("K31").Value = 1374.5625
$ws.Range("L31").Value = 4820.952
$ws.Range("M31").Value = -1079.5625
$ws.Range("N31").Value = -5410.952
$ws.Range("H34").Value = 4269.53
$ws.Range("I34").Value = 1374.5625
$ws.Range("J34").Value = 4820.952
$ws.Range("K34").Value = 1374.5625
$ws.Range("L34").Value = 4820.952
$ws.Range("M34").Value = -1172.5625
$ws.Range("N34").Value = -5224.952
$ws.Range("H134").Value = 3595.4285
$ws.Range("I134").Value = 3549.3572
$ws.Range("K134").Value = 10648.0716
$ws.Range("M134").Value = -8113.071599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 172.45454
$ws.Range("I23").Value = 115
$ws.Range("J23").Value = 185.22223
$ws.Range("K23").Value = 345
$ws.Range("L23").Value = 555.66669
$ws.Range("M23").Value = -110
$ws.Range("N23").Value = -1025.66669
$ws.Range("H56").Value = 41671784
$ws.Range("I56").Value = 41671784
$ws.Range("K56").Value = 41671784
$ws.Range("M56").Value = -41671254
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1883.4
$ws.Range("I93").Value = 1975.7858
$ws.Range("J93").Value = 590
$ws.Range("K93").Value = 1975.7858
$ws.Range("L93").Value = 590
$ws.Range("M93").Value = -727.7858000000001
$ws.Range("N93").Value = -3086
$ws.Range("H96").Value = 59994.5
$ws.Range("J96").Value = 59994.5
$ws.Range("L96").Value = 59994.5
$ws.Range("N96").Value = -65486.5
$ws.Range("H132").Value = 7841.375
$ws.Range("I132").Value = 11305.272
$ws.Range("K132").Value = 33915.81600000001
$ws.Range("M132").Value = -31385.81600000001
$ws.Range("H136").Value = 39257.43
$ws.Range("I136").Value = 58072.777
$ws.Range("K136").Value = 174218.331
$ws.Range("M136").Value = -171668.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 7250000
$ws.Range("I3").Value = 7250000
$ws.Range("K3").Value = 7250000
$ws.Range("M3").Value = -7249886
$ws.Range("H31").Value = 9999
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
